{"js": "// 1. Remove the \"Meta description\" paragraph that currently sits right\n//    after the title heading (paragraph index 1: \"Meta description: \u2026\").\nconst body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items/text\");\nawait context.sync();\n\nlet metaPara = null;\nfor (let i = 0; i < paras.items.length; i++) {\n  const t = paras.items[i].text;\n  if (t.indexOf(\"Meta description\") === 0) {\n    metaPara = paras.items[i];\n    break;\n  }\n}\nif (metaPara) {\n  metaPara.delete();\n  await context.sync();\n}\n\n// 2. At the end of the document, insert a new bold paragraph carrying\n//    the page title right before the final (italic) image-prompt\n//    paragraph, and swap that final paragraph's text for the old meta\n//    description copy (keeping its italic run formatting).\nconst paras2 = body.paragraphs;\nparas2.load(\"items/text\");\nawait context.sync();\n\nconst lastIndex = paras2.items.length - 1;\nconst lastPara = paras2.items[lastIndex];\nconst secondToLastPara = paras2.items[lastIndex - 1];\n\n// Insert the new bold paragraph right after the paragraph that currently\n// precedes the final paragraph \u2014 this lands it immediately before the\n// final paragraph without inheriting that paragraph's italic run format.\nconst insertionRange = secondToLastPara.getRange(Word.RangeLocation.end);\nconst newParaOoxml =\n  '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play FashionTV Highlife Slot Free \\u2013 Fashion and Style Slot Game</w:t></w:r></w:p>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>';\ninsertionRange.insertOoxml(newParaOoxml, Word.InsertLocation.after);\nawait context.sync();\n\n// Replace the final paragraph's text (was the \"Create a feature image\u2026\"\n// prompt) with the old meta-description copy, preserving its italic run.\nconst paras3 = body.paragraphs;\nparas3.load(\"items/text\");\nawait context.sync();\nconst finalPara = paras3.items[paras3.items.length - 1];\nconst finalRange = finalPara.getRange();\nfinalRange.insertText(\n  \"Experience exciting gameplay and special features on FashionTV Highlife Slot. Play free on desktop and mobile devices. Win with bonus features and special nudging function.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1. Remove the \"Meta description\" paragraph that currently sits right\n#    after the title heading.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.StartsWith(\"Meta description\")) {\n        $p.Range.Delete()\n        break\n    }\n}\n\n# 2. Insert a new bold paragraph carrying the page title right before the\n#    final (italic) image-prompt paragraph, without disturbing either of\n#    its neighboring paragraphs.\n$n = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($n)\n$startPos = $lastPara.Range.Start\n$insertionPoint = $d.Range($startPos, $startPos)\n\n$newParaOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n  '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n  '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n  '<pkg:xmlData>' +\n  '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' +\n  '<w:body>' +\n  '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play FashionTV Highlife Slot Free \u2013 Fashion and Style Slot Game</w:t></w:r></w:p>' +\n  '<w:p/>' +\n  '</w:body></w:document>' +\n  '</pkg:xmlData></pkg:part></pkg:package>'\n$insertionPoint.InsertXML($newParaOoxml)\n\n# InsertXML leaves a spare empty paragraph behind (its trailing paragraph\n# mark absorbs into the paragraph that used to start at the insertion\n# point) \u2014 drop that now-empty paragraph so the original final paragraph\n# is restored untouched right after our new paragraph.\n$spacerIndex = $n + 1\n$spacerPara = $d.Paragraphs.Item($spacerIndex)\nif ($spacerPara.Range.Text.Trim().Length -eq 0) {\n    $spacerPara.Range.Delete()\n}\n\n# 3. Swap the final paragraph's text (previously the \"Create a feature\n#    image\u2026\" prompt) for the old meta-description copy, while preserving\n#    its italic run formatting via Find/Replace scoped to that paragraph.\n$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$find = $finalPara.Range.Find\n$findText = \"Create a feature image for FashionTV Highlife slot game that captures the luxurious and glamorous atmosphere of the game. The image should be in a cartoon style, depicting a happy Maya warrior with glasses, surrounded by the symbols of the game, such as a car, yacht, gold watch, and the three fantastic models. The background should feature the FashionTV Highlife logo and a vibrant city at night with bright lights shining. Make sure to include the Special Nudging HP1 and Free Spins function icons in the image, emphasizing the special features of the game. Overall, the feature image should convey the excitement of winning high sums while enjoying the extravagance and luxury of this online slot game.\"\n$replaceText = \"Experience exciting gameplay and special features on FashionTV Highlife Slot. Play free on desktop and mobile devices. Win with bonus features and special nudging function.\"\n$find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n"}
